# Updated cryptos list values (price + 1h volume-change%) to match the
# latest scrape. A handful of rows also had their Coin/Link swapped
# (Chainlink <-> Dai moved between rows 19 and 20).
#
# Price cells that look like plain numbers (e.g. "7.20", "0.0620") are
# written with a leading apostrophe so Excel keeps them as text and
# preserves trailing zeros, matching how the sheet already stores them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.740.00'
$ws.Range('E2').Value = '  +0.28%  '

# Row 3
$ws.Range('D3').Value = '1.602.48'
$ws.Range('E3').Value = '  +0.22%  '

# Row 4
$ws.Range('E4').Value = '  +0.19%  '

# Row 5
$ws.Range('D5').Value = '''211.83'
$ws.Range('E5').Value = '  +0.00%  '

# Row 6
$ws.Range('D6').Value = '''0.513'
$ws.Range('E6').Value = '  -0.45%  '

# Row 7
$ws.Range('E7').Value = '  +0.18%  '

# Row 8
$ws.Range('D8').Value = '''0.0620'
$ws.Range('E8').Value = '  +0.18%  '

# Row 10
$ws.Range('D10').Value = '''19.76'
$ws.Range('E10').Value = '  +0.92%  '

# Row 11
$ws.Range('D11').Value = '''0.0846'
$ws.Range('E11').Value = '  +0.85%  '

# Row 12
$ws.Range('D12').Value = '1.827.10'
$ws.Range('E12').Value = '  +0.19%  '

# Row 13
$ws.Range('D13').Value = '1.599.87'
$ws.Range('E13').Value = '  +0.23%  '

# Row 14
$ws.Range('D14').Value = '''4.05'
$ws.Range('E14').Value = '  +0.46%  '

# Row 15
$ws.Range('E15').Value = '  -0.35%  '

# Row 16
$ws.Range('D16').Value = '''65.09'
$ws.Range('E16').Value = '  -0.09%  '

# Row 17
$ws.Range('E17').Value = '  +0.41%  '

# Row 18
$ws.Range('D18').Value = '''210.37'
$ws.Range('E18').Value = '  +0.73%  '

# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''7.20'
$ws.Range('E19').Value = '  +2.31%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.01'
$ws.Range('E20').Value = '  +0.21%  '

# Row 21
$ws.Range('D21').Value = '''4.29'
$ws.Range('E21').Value = '  -0.01%  '

# Row 22
$ws.Range('D22').Value = '''2.28'
$ws.Range('E22').Value = '  -1.76%  '

# Row 23
$ws.Range('E23').Value = '  +0.27%  '

# Row 24
$ws.Range('D24').Value = '''143.65'
$ws.Range('E24').Value = '  -1.14%  '

# Row 25
$ws.Range('E25').Value = '  +0.10%  '

# Row 26
$ws.Range('E26').Value = '  -0.16%  '

# Row 27
$ws.Range('E27').Value = '  -0.83%  '

# Row 28
$ws.Range('D28').Value = '''15.40'
$ws.Range('E28').Value = '  +0.71%  '

# Row 29
$ws.Range('D29').Value = '''0.0512'
$ws.Range('E29').Value = '  -0.35%  '

# Row 30
$ws.Range('E30').Value = '  +0.04%  '

# Row 31
$ws.Range('E31').Value = '  +1.34%  '

# Row 32
$ws.Range('E32').Value = '  +1.18%  '

# Row 33
$ws.Range('D33').Value = '1.295.78'
$ws.Range('E33').Value = '  +1.51%  '

# Row 34
$ws.Range('E34').Value = '  +0.78%  '

# Row 35
$ws.Range('E35').Value = '  +0.87%  '

# Row 36
$ws.Range('D36').Value = '''0.602'
$ws.Range('E36').Value = '  -2.96%  '

# Row 37
$ws.Range('E37').Value = '  +11.07%  '

# Row 38
$ws.Range('E38').Value = '  -0.22%  '

# Row 39
$ws.Range('D39').Value = '''0.832'
$ws.Range('E39').Value = '  -0.33%  '

# Row 40
$ws.Range('E40').Value = '  -1.73%  '

# Row 41
$ws.Range('E41').Value = '  -0.21%  '

# Row 42
$ws.Range('D42').Value = '''0.786'
$ws.Range('E42').Value = '  +0.10%  '

# Row 43
$ws.Range('D43').Value = '''63.11'
$ws.Range('E43').Value = '  -1.35%  '

# Row 44
$ws.Range('D44').Value = '1.739.17'
$ws.Range('E44').Value = '  +0.17%  '

# Row 45
$ws.Range('D45').Value = '''90.60'
$ws.Range('E45').Value = '  -0.71%  '

# Row 46
$ws.Range('E46').Value = '  -2.11%  '

# Row 47
$ws.Range('E47').Value = '  -0.35%  '

# Row 48
$ws.Range('D48').Value = '''0.0517'
$ws.Range('E48').Value = '  +1.87%  '

# Row 49
$ws.Range('E49').Value = '  +0.09%  '

# Row 50
$ws.Range('D50').Value = '''7.45'
$ws.Range('E50').Value = '  +0.68%  '

# Row 51
$ws.Range('E51').Value = '  +0.83%  '
